# Append 6 new arrival rows (325..330, 1-based "NUMBER" 325-330) to the
# "Main Data" sheet of the KTW_Arrivals workbook, continuing directly after
# the existing last row (Excel row 325).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$date = "Monday, Jan 16"

# Each entry: Number, Time, Flight, From, Short, Airline, Model, AircraftID, Status, Difference
$newRows = @(
    @(325, "1:20 AM", "QY5916", "Leipzig",  "(LEJ)", "Bluebird Nordic ", "B738", "(TF-BBQ)", "1:00 AM", "0 hours, -20 minutes"),
    @(326, "2:57 AM", "E45160", "Hurghada", "(HRG)", "Enter Air ",       "B738", "(SP-ESH)", "2:07 AM", "0 hours, -50 minutes"),
    @(327, "3:20 AM", "BO624",  "Catania",  "(CTA)", "Bluebird Nordic ", "B734", "(TF-BBN)", "2:47 AM", "0 hours, -33 minutes"),
    @(328, "6:00 AM", "SAR1980","Cologne",  "(CGN)", "SprintAir ",       "AT75", "(SP-SPF)", "5:41 AM", "0 hours, -19 minutes"),
    @(329, "6:31 AM", "UNKNOWN","Budapest", "(BUD)", "Wizz Air ",        "A21N", "(HA-LZD)", "6:04 AM", "0 hours, -27 minutes"),
    @(330, "8:19 AM", "BO950",  "Cagliari", "(CAG)", "Bluebird Nordic ", "B734", "(TF-BBO)", "8:07 AM", "0 hours, -12 minutes")
)

$startRow = 326

for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $data = $newRows[$i]

    $ws.Cells.Item($r, 1).Value = $data[0]
    $ws.Cells.Item($r, 2).Value = $date
    $ws.Cells.Item($r, 3).Value = $data[1]
    $ws.Cells.Item($r, 4).Value = $data[2]
    $ws.Cells.Item($r, 5).Value = $data[3]
    $ws.Cells.Item($r, 6).Value = $data[4]
    $ws.Cells.Item($r, 7).Value = $data[5]
    $ws.Cells.Item($r, 8).Value = $data[6]
    $ws.Cells.Item($r, 9).Value = $data[7]
    $ws.Cells.Item($r, 10).Value = $data[8]
    $ws.Cells.Item($r, 12).Value = $data[9]
}
